$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 60: a new time log entry
$ws.Range("B60").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C60").Value = 0.5
$ws.Range("D60").Value = "Finished 1 small problem"

# Update the active selection to match the saved view (D60)
$ws.Range("D60").Select()
